$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Gnai2"
$ws.Cells.Item(2,3).Value = "Agtr2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 195.0792385
$ws.Cells.Item(2,8).Value = 390.158477
$ws.Cells.Item(2,9).Value = 0.2640605522989327
$ws.Cells.Item(2,10).Value = 0.1982306263353075
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.5
$ws.Cells.Item(2,13).Value = 0.1409635
$ws.Cells.Item(2,14).Value = 0.281927
$ws.Cells.Item(2,15).Value = 0.127302058387171
$ws.Cells.Item(2,16).Value = 0.0886289202516707
$ws.Cells.Item(2,17).Value = 27.49905223629475
$ws.Cells.Item(2,18).Value = 109.996208945179
$ws.Cells.Item(2,19).Value = 0.03361545184650734
$ws.Cells.Item(2,20).Value = 0.0175689663729107
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Gnai2"
$ws.Cells.Item(3,3).Value = "Agtr2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 195.0792385
$ws.Cells.Item(3,8).Value = 390.158477
$ws.Cells.Item(3,9).Value = 0.2640605522989327
$ws.Cells.Item(3,10).Value = 0.1982306263353075
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.9663516666666667
$ws.Cells.Item(3,14).Value = 2.899055
$ws.Cells.Item(3,15).Value = 0.872697941612829
$ws.Cells.Item(3,16).Value = 0.9113710797483293
$ws.Cells.Item(3,17).Value = 188.5151472565392
$ws.Cells.Item(3,18).Value = 1131.090883539235
$ws.Cells.Item(3,19).Value = 0.2304451004524253
$ws.Cells.Item(3,20).Value = 0.1806616599623968
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Gnai2"
$ws.Cells.Item(4,3).Value = "Agtr2"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 62.40792233333334
$ws.Cells.Item(4,8).Value = 187.223767
$ws.Cells.Item(4,9).Value = 0.08447577797556809
$ws.Cells.Item(4,10).Value = 0.09512412720758515
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.5
$ws.Cells.Item(4,13).Value = 0.1409635
$ws.Cells.Item(4,14).Value = 0.281927
$ws.Cells.Item(4,15).Value = 0.127302058387171
$ws.Cells.Item(4,16).Value = 0.0886289202516707
$ws.Cells.Item(4,17).Value = 8.797239159834833
$ws.Cells.Item(4,18).Value = 52.783434959009
$ws.Cells.Item(4,19).Value = 0.01075394042014746
$ws.Cells.Item(4,20).Value = 0.008430748684290843
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Gnai2"
$ws.Cells.Item(5,3).Value = "Agtr2"
$ws.Cells.Item(5,4).Value = "FAPs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 62.40792233333334
$ws.Cells.Item(5,8).Value = 187.223767
$ws.Cells.Item(5,9).Value = 0.08447577797556809
$ws.Cells.Item(5,10).Value = 0.09512412720758515
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.9663516666666667
$ws.Cells.Item(5,14).Value = 2.899055
$ws.Cells.Item(5,15).Value = 0.872697941612829
$ws.Cells.Item(5,16).Value = 0.9113710797483293
$ws.Cells.Item(5,17).Value = 60.30799976002056
$ws.Cells.Item(5,18).Value = 542.771997840185
$ws.Cells.Item(5,19).Value = 0.07372183755542062
$ws.Cells.Item(5,20).Value = 0.0866933785232943
$ws.Cells.Item(6,1).Value = "M1"
$ws.Cells.Item(6,2).Value = "Gnai2"
$ws.Cells.Item(6,3).Value = "Agtr2"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 142.8621113333333
$ws.Cells.Item(6,8).Value = 428.586334
$ws.Cells.Item(6,9).Value = 0.1933791023142199
$ws.Cells.Item(6,10).Value = 0.2177549443006804
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.5
$ws.Cells.Item(6,13).Value = 0.1409635
$ws.Cells.Item(6,14).Value = 0.281927
$ws.Cells.Item(6,15).Value = 0.127302058387171
$ws.Cells.Item(6,16).Value = 0.0886289202516707
$ws.Cells.Item(6,17).Value = 20.13834323093633
$ws.Cells.Item(6,18).Value = 120.830059385618
$ws.Cells.Item(6,19).Value = 0.02461755777366353
$ws.Cells.Item(6,20).Value = 0.019299385592832
$ws.Cells.Item(7,1).Value = "M1"
$ws.Cells.Item(7,2).Value = "Gnai2"
$ws.Cells.Item(7,3).Value = "Agtr2"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 142.8621113333333
$ws.Cells.Item(7,8).Value = 428.586334
$ws.Cells.Item(7,9).Value = 0.1933791023142199
$ws.Cells.Item(7,10).Value = 0.2177549443006804
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.9663516666666667
$ws.Cells.Item(7,14).Value = 2.899055
$ws.Cells.Item(7,15).Value = 0.872697941612829
$ws.Cells.Item(7,16).Value = 0.9113710797483293
$ws.Cells.Item(7,17).Value = 138.0550393904856
$ws.Cells.Item(7,18).Value = 1242.49535451437
$ws.Cells.Item(7,19).Value = 0.1687615445405564
$ws.Cells.Item(7,20).Value = 0.1984555587078484
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Gnai2"
$ws.Cells.Item(8,3).Value = "Agtr2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 141.6168416666667
$ws.Cells.Item(8,8).Value = 424.850525
$ws.Cells.Item(8,9).Value = 0.1916934970264942
$ws.Cells.Item(8,10).Value = 0.2158568649262854
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.5
$ws.Cells.Item(8,13).Value = 0.1409635
$ws.Cells.Item(8,14).Value = 0.281927
$ws.Cells.Item(8,15).Value = 0.127302058387171
$ws.Cells.Item(8,16).Value = 0.0886289202516707
$ws.Cells.Item(8,17).Value = 19.96280566027917
$ws.Cells.Item(8,18).Value = 119.776833961675
$ws.Cells.Item(8,19).Value = 0.02440297675090775
$ws.Cells.Item(8,20).Value = 0.0191311608673274
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Gnai2"
$ws.Cells.Item(9,3).Value = "Agtr2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 141.6168416666667
$ws.Cells.Item(9,8).Value = 424.850525
$ws.Cells.Item(9,9).Value = 0.1916934970264942
$ws.Cells.Item(9,10).Value = 0.2158568649262854
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.9663516666666667
$ws.Cells.Item(9,14).Value = 2.899055
$ws.Cells.Item(9,15).Value = 0.872697941612829
$ws.Cells.Item(9,16).Value = 0.9113710797483293
$ws.Cells.Item(9,17).Value = 136.8516709726528
$ws.Cells.Item(9,18).Value = 1231.665038753875
$ws.Cells.Item(9,19).Value = 0.1672905202755864
$ws.Cells.Item(9,20).Value = 0.196725704058958
$ws.Cells.Item(10,1).Value = "Neutro"
$ws.Cells.Item(10,2).Value = "Gnai2"
$ws.Cells.Item(10,3).Value = "Agtr2"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 143.783834
$ws.Cells.Item(10,8).Value = 431.351502
$ws.Cells.Item(10,9).Value = 0.1946267522348261
$ws.Cells.Item(10,10).Value = 0.2191598631141254
$ws.Cells.Item(10,11).Value = 1
$ws.Cells.Item(10,12).Value = 0.5
$ws.Cells.Item(10,13).Value = 0.1409635
$ws.Cells.Item(10,14).Value = 0.281927
$ws.Cells.Item(10,15).Value = 0.127302058387171
$ws.Cells.Item(10,16).Value = 0.0886289202516707
$ws.Cells.Item(10,17).Value = 20.268272484059
$ws.Cells.Item(10,18).Value = 121.609634904354
$ws.Cells.Item(10,19).Value = 0.02477638617670329
$ws.Cells.Item(10,20).Value = 0.01942390203030889
$ws.Cells.Item(11,1).Value = "Neutro"
$ws.Cells.Item(11,2).Value = "Gnai2"
$ws.Cells.Item(11,3).Value = "Agtr2"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 143.783834
$ws.Cells.Item(11,8).Value = 431.351502
$ws.Cells.Item(11,9).Value = 0.1946267522348261
$ws.Cells.Item(11,10).Value = 0.2191598631141254
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.9663516666666667
$ws.Cells.Item(11,14).Value = 2.899055
$ws.Cells.Item(11,15).Value = 0.872697941612829
$ws.Cells.Item(11,16).Value = 0.9113710797483293
$ws.Cells.Item(11,17).Value = 138.9457476256233
$ws.Cells.Item(11,18).Value = 1250.51172863061
$ws.Cells.Item(11,19).Value = 0.1698503660581228
$ws.Cells.Item(11,20).Value = 0.1997359610838165
$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Gnai2"
$ws.Cells.Item(12,3).Value = "Agtr2"
$ws.Cells.Item(12,4).Value = "ECs"
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 53.01711450000001
$ws.Cells.Item(12,8).Value = 106.034229
$ws.Cells.Item(12,9).Value = 0.07176431814995911
$ws.Cells.Item(12,10).Value = 0.05387357411601602
$ws.Cells.Item(12,11).Value = 1
$ws.Cells.Item(12,12).Value = 0.5
$ws.Cells.Item(12,13).Value = 0.1409635
$ws.Cells.Item(12,14).Value = 0.281927
$ws.Cells.Item(12,15).Value = 0.127302058387171
$ws.Cells.Item(12,16).Value = 0.0886289202516707
$ws.Cells.Item(12,17).Value = 7.47347801982075
$ws.Cells.Item(12,18).Value = 29.893912079283
$ws.Cells.Item(12,19).Value = 0.009135745419241609
$ws.Cells.Item(12,20).Value = 0.004774756704000854
$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Gnai2"
$ws.Cells.Item(13,3).Value = "Agtr2"
$ws.Cells.Item(13,4).Value = "FAPs"
$ws.Cells.Item(13,5).Value = 2
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 53.01711450000001
$ws.Cells.Item(13,8).Value = 106.034229
$ws.Cells.Item(13,9).Value = 0.07176431814995911
$ws.Cells.Item(13,10).Value = 0.05387357411601602
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.9663516666666667
$ws.Cells.Item(13,14).Value = 2.899055
$ws.Cells.Item(13,15).Value = 0.872697941612829
$ws.Cells.Item(13,16).Value = 0.9113710797483293
$ws.Cells.Item(13,17).Value = 51.23317695893251
$ws.Cells.Item(13,18).Value = 307.399061753595
$ws.Cells.Item(13,19).Value = 0.06262857273071749
$ws.Cells.Item(13,20).Value = 0.04909881741201516